$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The accession number / DOI for the reference genome row (row 2, column D)
# was updated from the old GCA accession to the new one.
$ws.Range("D2").Value = "GCA_963668995.1"

# Reflect the active cell selection recorded in the saved file.
$ws.Range("D2").Select()
